$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("展览")
$ws.Range("F5").Value = 744
$ws.Range("F6").Value = 2332
$ws.Range("F8").Value = 1736
$ws.Range("F9").Value = 2929
$ws.Range("F10").Value = 167
$ws.Range("F11").Value = 4385
$ws.Range("F12").Value = 378
$ws.Range("F13").Value = 208
$ws.Range("F14").Value = 139
$ws.Range("F15").Value = 553
$ws.Range("F16").Value = 262
$ws.Range("F17").Value = 13
$ws.Range("F18").Value = 182
$ws.Range("F20").Value = 105
$ws.Range("F22").Value = 4398
$ws.Range("F24").Value = 3650
$ws.Range("F25").Value = 1133
$ws.Range("F26").Value = 211
$ws.Range("F27").Value = 558
$ws.Range("F29").Value = 86
$ws.Range("F30").Value = 551
$ws.Range("F31").Value = 552
$ws.Range("F32").Value = 511

$ws = $wb.Worksheets.Item("演出")
$ws.Range("F2").Value = 4

$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F3").Value = 1035
$ws.Range("F4").Value = 15

$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F4").Value = 1035
$ws.Range("F5").Value = 15
$ws.Range("F8").Value = 744
$ws.Range("F9").Value = 2332
$ws.Range("F11").Value = 1736
$ws.Range("F12").Value = 4
$ws.Range("F13").Value = 2929
$ws.Range("F14").Value = 167
$ws.Range("F15").Value = 4385
$ws.Range("F16").Value = 378
$ws.Range("F17").Value = 208
$ws.Range("F18").Value = 139
$ws.Range("F19").Value = 553
$ws.Range("F20").Value = 262
$ws.Range("F21").Value = 13
$ws.Range("F22").Value = 182
$ws.Range("F25").Value = 105
$ws.Range("F27").Value = 4398
$ws.Range("F29").Value = 3650
$ws.Range("F30").Value = 1133
$ws.Range("F31").Value = 211
$ws.Range("F32").Value = 558
$ws.Range("F34").Value = 86
$ws.Range("F35").Value = 551
$ws.Range("F36").Value = 552
$ws.Range("F37").Value = 511
